$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Styles: add a new yellow fill + cell style (fills 3->4, cellXfs 3->4) ---
# This will be realized implicitly by applying a new Interior color to E2 below,
# which causes the engine to allocate a new fill + cellXf entry.

# --- Column width changes ---
$ws.Columns.Item(3).ColumnWidth = 80.1666666666667   # C: 56 -> 81
$ws.Columns.Item(4).ColumnWidth = 37.1666666666667   # D: 67 -> 38
$ws.Columns.Item(6).ColumnWidth = 16.1666666666667   # F: 16 -> 17
$ws.Columns.Item(8).ColumnWidth = 44.1666666666667   # H: 60 -> 45

# --- Row data updates (rows 2-13) ---

$ws.Range("A2").Value = "1328395"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328395"
$ws.Range("C2").Value = "Professional European Key Account Desk Expert (EU Preferred)"
$ws.Range("D2").Value = "Maastricht, Netherlands"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "12 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "DHL Group"

$ws.Range("A3").Value = "1327778"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1327778"
$ws.Range("C3").Value = "Digital Content & Stakeholder Engagement Intern"
$ws.Range("D3").Value = "Colombo, Sri Lanka"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "12 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "Solutions Ground (Pvt) Ltd"

$ws.Range("A4").Value = "1327475"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1327475"
$ws.Range("C4").Value = "Property Consultant"
$ws.Range("D4").Value = "Cairo, Cairo Governorate, Egypt"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "10 applicants"
$ws.Range("G4").Value = "9 - 12 Weeks"
$ws.Range("H4").Value = "Bold Routes"

$ws.Range("A5").Value = "1327366"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1327366"
$ws.Range("C5").Value = "Marketing with advertisement integrated maps and practical experience in Japan"
$ws.Range("D5").Value = "日本、愛知県名古屋市"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "55 applicants"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "HYOJITO Co.,Ltd."

$ws.Range("A6").Value = "1325908"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1325908"
$ws.Range("C6").Value = "International Business Development Intern"
$ws.Range("D6").Value = "Colombo, Sri Lanka"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "18 applicants"
$ws.Range("G6").Value = "9 - 12 Weeks"
$ws.Range("H6").Value = "Fintechnology Asia Pacific Lanka (Pvt) Ltd"

$ws.Range("A7").Value = "1325856"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1325856"
$ws.Range("C7").Value = "SALES ASSISTANT"
$ws.Range("D7").Value = "Denizli, Kumkısık, Denizli, Türkiye"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "53 applicants"
$ws.Range("G7").Value = "6 - 18 Months"
$ws.Range("H7").Value = "MULBERRY HOME"

$ws.Range("A8").Value = "1325702"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1325702"
$ws.Range("C8").Value = "Guest Relations Executive and Waitress"
$ws.Range("D8").Value = "Colombo, Sri Lanka"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "10 applicants"
$ws.Range("G8").Value = "3 - 6 Months"
$ws.Range("H8").Value = "Indian Kitchen PVT LTD"

$ws.Range("A9").Value = "1325417"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1325417"
$ws.Range("C9").Value = "Junior Software Engineer – AI & Internal Tools (EU ONLY)"
$ws.Range("D9").Value = "Brussels, Belgium"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "126 applicants"
$ws.Range("G9").Value = "9 - 12 Weeks"
$ws.Range("H9").Value = "Eureka Resource Mining"

$ws.Range("A10").Value = "1323735"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1323735"
$ws.Range("C10").Value = "Export & Sales Support Assistant"
$ws.Range("D10").Value = "İstanbul, Türkiye"
$ws.Range("E10").Value = "No"
$ws.Range("F10").Value = "131 applicants"
$ws.Range("G10").Value = "9 - 12 Weeks"
$ws.Range("H10").Value = "Anıl Lingerie & Homewear"

$ws.Range("A11").Value = "1321052"
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1321052"
$ws.Range("C11").Value = "International Sales Representetive"
$ws.Range("D11").Value = "Maslak, Sarıyer/İstanbul, Türkiye"
$ws.Range("E11").Value = "No"
$ws.Range("F11").Value = "128 applicants"
$ws.Range("G11").Value = "6 - 18 Months"
$ws.Range("H11").Value = "Esvita Clinic"

$ws.Range("A12").Value = "1289380"
$ws.Range("B12").Value = "https://aiesec.org/opportunity/global-talent/1289380"
$ws.Range("C12").Value = "Medical Advisor Polish Speaker"
$ws.Range("D12").Value = "İstanbul, Türkiye"
$ws.Range("E12").Value = "No"
$ws.Range("F12").Value = "6 applicants"
$ws.Range("G12").Value = "6 - 18 Months"
$ws.Range("H12").Value = "International Plus"

$ws.Range("A13").Value = "1289375"
$ws.Range("B13").Value = "https://aiesec.org/opportunity/global-talent/1289375"
$ws.Range("C13").Value = "Medical Advisor (German Speaker)"
$ws.Range("D13").Value = "İstanbul, Türkiye"
$ws.Range("E13").Value = "No"
$ws.Range("F13").Value = "28 applicants"
$ws.Range("G13").Value = "6 - 18 Months"
$ws.Range("H13").Value = "International Plus"

# --- E2 highlight: new yellow fill style applied only to E2 ---
$ws.Range("E2").Interior.ColorIndex = 6   # yellow (RGB 255,255,0)

Write-Host "Edit complete"
